$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 ("Age" section header) and Row 28 ("Education" section header), three
# languages each, get "By ..." wording. Edit order matches the shared-string
# append order seen in the target workbook.
$ws.Range("C18").Value = "By age (in years) "
$ws.Range("C28").Value = "By education"
$ws.Range("B18").Value = "По возрасту (в годах)"
$ws.Range("A18").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A28").Value = "Билими боюнча"
$ws.Range("B28").Value = "По образованию"
